$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.615.60"
$ws.Range("E2").Value = "  +3.62%  "

$ws.Range("D3").Value = "3.075.35"
$ws.Range("E3").Value = "  +4.51%  "

$ws.Range("E4").Value = "  -0.67%  "

$ws.Range("D5").Value = "517.39"
$ws.Range("E5").Value = "  +3.77%  "

$ws.Range("D6").Value = "141.02"
$ws.Range("E6").Value = "  +3.60%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.54%  "

$ws.Range("E8").Value = "  +1.51%  "

$ws.Range("D9").Value = "7.28"
$ws.Range("E9").Value = "  +4.11%  "

$ws.Range("E10").Value = "  +2.84%  "

$ws.Range("D11").Value = "0.374"
$ws.Range("E11").Value = "  +3.19%  "

$ws.Range("D12").Value = "3.605.27"
$ws.Range("E12").Value = "  +2.51%  "

$ws.Range("E13").Value = "  +3.19%  "

$ws.Range("D14").Value = "25.66"
$ws.Range("E14").Value = "  -0.95%  "

$ws.Range("E15").Value = "  +2.55%  "

$ws.Range("D16").Value = "57.675.60"
$ws.Range("E16").Value = "  +2.91%  "

$ws.Range("D17").Value = "3.077.22"
$ws.Range("E17").Value = "  +2.14%  "

$ws.Range("D18").Value = "6.07"
$ws.Range("E18").Value = "  +2.98%  "

$ws.Range("D19").Value = "13.00"
$ws.Range("E19").Value = "  +1.24%  "

$ws.Range("D20").Value = "8.09"
$ws.Range("E20").Value = "  +4.17%  "

$ws.Range("D21").Value = "334.51"
$ws.Range("E21").Value = "  +3.67%  "

$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.43%  "

$ws.Range("E23").Value = "  +2.94%  "

$ws.Range("D24").Value = "65.95"
$ws.Range("E24").Value = "  +2.59%  "

$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  +4.85%  "

$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -1.12%  "

$ws.Range("D27").Value = "0.0₃0911"
$ws.Range("E27").Value = "  +5.10%  "

$ws.Range("E28").Value = "  -0.76%  "

$ws.Range("E29").Value = "  +3.94%  "

$ws.Range("E30").Value = "  +3.79%  "

$ws.Range("D31").Value = "20.88"
$ws.Range("E31").Value = "  +4.40%  "

$ws.Range("E32").Value = "  +2.38%  "

$ws.Range("D33").Value = "154.66"
$ws.Range("E33").Value = "  +2.09%  "

$ws.Range("B34").Value = "EnergySwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D34").Value = "27.17"
$ws.Range("E34").Value = "  +8.99%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "4.47"
$ws.Range("E35").Value = "  -1.06%  "

$ws.Range("D36").Value = "5.90"
$ws.Range("E36").Value = "  +3.42%  "

$ws.Range("E37").Value = "  +3.84%  "

$ws.Range("E38").Value = "  +3.75%  "

$ws.Range("D39").Value = "3.115.44"
$ws.Range("E39").Value = "  +3.34%  "

$ws.Range("D40").Value = "3.92"
$ws.Range("E40").Value = "  +5.53%  "

$ws.Range("D41").Value = "36.97"
$ws.Range("E41").Value = "  +1.35%  "

$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.91%  "

$ws.Range("D43").Value = "0.659"
$ws.Range("E43").Value = "  +1.36%  "

$ws.Range("D44").Value = "2.264.96"
$ws.Range("E44").Value = "  +5.39%  "

$ws.Range("E45").Value = "  +8.23%  "

$ws.Range("E46").Value = "  +3.37%  "

$ws.Range("D47").Value = "20.04"
$ws.Range("E47").Value = "  +2.82%  "

$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "5.87"
$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "0.927"
$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("D50").Value = "265.43"
$ws.Range("E50").Value = "  +17.68%  "

$ws.Range("E51").Value = "  +3.60%  "
